# Fill in the simulated-game transition-matrix probabilities.
# Only the cells that move off their initial 0 value are touched; every
# other cell in the matrix keeps its existing value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Starting_State = Af0)
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.4666666666666667
$ws.Range("P2").Value = 0.1333333333333333
$ws.Range("S2").Value = 0.2

# Row 3 (Starting_State = Af1)
$ws.Range("C3").Value = 0.125
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.375

# Row 4 (Starting_State = Af2)
$ws.Range("P4").Value = 1

# Row 6 (Starting_State = Ai0)
$ws.Range("B6").Value = 0.1333333333333333
$ws.Range("J6").Value = 0.5333333333333333
$ws.Range("R6").Value = 0.1333333333333333
$ws.Range("S6").Value = 0.2

# Row 7 (Starting_State = Ai1)
$ws.Range("F7").Value = 0.09090909090909091
$ws.Range("J7").Value = 0.2727272727272727
$ws.Range("S7").Value = 0.6363636363636364

# Row 8 (Starting_State = Ai2)
$ws.Range("B8").Value = 0.05084745762711865
$ws.Range("D8").Value = 0.01694915254237288
$ws.Range("F8").Value = 0.03389830508474576
$ws.Range("J8").Value = 0.1016949152542373
$ws.Range("Q8").Value = 0.1186440677966102
$ws.Range("R8").Value = 0.2033898305084746
$ws.Range("S8").Value = 0.4745762711864407

# Row 9 (Starting_State = Ai3)
$ws.Range("B9").Value = 0.08333333333333333
$ws.Range("F9").Value = 0.08333333333333333
$ws.Range("J9").Value = 0.08333333333333333
$ws.Range("S9").Value = 0.75

# Row 10 (Starting_State = Ar0)
$ws.Range("B10").Value = 0.07317073170731707
$ws.Range("F10").Value = 0.08536585365853659
$ws.Range("J10").Value = 0.1463414634146341
$ws.Range("O10").Value = 0.02439024390243903
$ws.Range("Q10").Value = 0.1463414634146341
$ws.Range("R10").Value = 0.0975609756097561
$ws.Range("S10").Value = 0.4268292682926829

# Row 11 (Starting_State = Bf0)
$ws.Range("G11").Value = 0.1764705882352941
$ws.Range("K11").Value = 0.1764705882352941
$ws.Range("L11").Value = 0.6470588235294118

# Row 12 (Starting_State = Bf1)
$ws.Range("G12").Value = 0.7272727272727273
$ws.Range("J12").Value = 0.1818181818181818
$ws.Range("S12").Value = 0.09090909090909091

# Row 13 (Starting_State = Bf2)
$ws.Range("G13").Value = 0.25
$ws.Range("J13").Value = 0.75

# Row 15 (Starting_State = Bi0)
$ws.Range("H15").Value = 0.4285714285714285
$ws.Range("J15").Value = 0.3571428571428572
$ws.Range("S15").Value = 0.2142857142857143

# Row 16 (Starting_State = Bi1)
$ws.Range("F16").Value = 0.1428571428571428
$ws.Range("H16").Value = 0.2857142857142857
$ws.Range("J16").Value = 0.2857142857142857
$ws.Range("K16").Value = 0.1428571428571428
$ws.Range("O16").Value = 0.1428571428571428

# Row 17 (Starting_State = Bi2)
$ws.Range("H17").Value = 0.2105263157894737
$ws.Range("I17").Value = 0.1052631578947368
$ws.Range("J17").Value = 0.2105263157894737
$ws.Range("K17").Value = 0.1052631578947368
$ws.Range("M17").Value = 0.05263157894736842
$ws.Range("O17").Value = 0.1052631578947368
$ws.Range("S17").Value = 0.2105263157894737

# Row 18 (Starting_State = Bi3)
$ws.Range("H18").Value = 0.3636363636363636
$ws.Range("I18").Value = 0.1363636363636364
$ws.Range("J18").Value = 0.2272727272727273
$ws.Range("K18").Value = 0.04545454545454546
$ws.Range("O18").Value = 0.04545454545454546
$ws.Range("S18").Value = 0.1818181818181818

# Row 19 (Starting_State = Br0)
$ws.Range("F19").Value = 0.01785714285714286
$ws.Range("H19").Value = 0.3482142857142857
$ws.Range("I19").Value = 0.07142857142857142
$ws.Range("J19").Value = 0.2946428571428572
$ws.Range("K19").Value = 0.08928571428571429
$ws.Range("M19").Value = 0.01785714285714286
$ws.Range("O19").Value = 0.04464285714285714
$ws.Range("S19").Value = 0.1160714285714286
